$wb = $excel.ActiveWorkbook

# Sheet 1: "HSV Log OLS"
$ws1 = $wb.Worksheets.Item("HSV Log OLS")
$ws1.Range("A2").Value = 0.152
$ws1.Range("B2").Value = 4.716
$ws1.Range("C2").Value = 0.787

# Sheet 2: "HSV PPML"
$ws2 = $wb.Worksheets.Item("HSV PPML")
$ws2.Range("A2").Value = 0.024
$ws2.Range("B2").Value = 1.036
$ws2.Range("C2").Value = 0.781

# Sheet 3: "HSVT NLLSQ"
$ws3 = $wb.Worksheets.Item("HSVT NLLSQ")
$ws3.Range("A2").Value = -0.078
$ws3.Range("B2").Value = 0.244
$ws3.Range("C2").Value = 19341.44
$ws3.Range("D2").Value = 0.148
